$d = $word.ActiveDocument

# Remove the redundant lowercase "constel·lació," (with trailing comma) that
# duplicated the following "Constel·lació" word, leaving the surrounding
# spaces intact (resulting in a double space).
$d.Content.Find.Execute("constel·lació,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
